$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.722.40"
$ws.Range("E2").Value = "  -3.92%  "
$ws.Range("D3").Value = "2.905.22"
$ws.Range("E3").Value = "  -4.26%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.15"
$ws.Range("E5").Value = "  -0.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.28"
$ws.Range("E6").Value = "  -6.56%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -2.56%  "
$ws.Range("D9").Value = "2.903.88"
$ws.Range("E9").Value = "  -4.21%  "
$ws.Range("E10").Value = "  -4.70%  "
$ws.Range("E11").Value = "  -5.08%  "
$ws.Range("E12").Value = "  -4.62%  "
$ws.Range("E13").Value = "  -4.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.40"
$ws.Range("E14").Value = "  -6.54%  "
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("D16").Value = "3.383.19"
$ws.Range("E16").Value = "  -4.33%  "
$ws.Range("D17").Value = "60.684.57"
$ws.Range("E17").Value = "  -3.89%  "
$ws.Range("E18").Value = "  -6.02%  "
$ws.Range("D19").Value = "2.903.50"
$ws.Range("E19").Value = "  -4.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "427.43"
$ws.Range("E20").Value = "  -5.81%  "
$ws.Range("E21").Value = "  -5.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.684"
$ws.Range("E22").Value = "  -2.12%  "
$ws.Range("E23").Value = "  -6.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.18"
$ws.Range("E24").Value = "  -2.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.78"
$ws.Range("E25").Value = "  -6.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.21"
$ws.Range("E26").Value = "  -6.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.89"
$ws.Range("E27").Value = "  -4.32%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  -2.86%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  -3.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.05"
$ws.Range("E32").Value = "  -7.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.40"
$ws.Range("E33").Value = "  -4.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.107"
$ws.Range("E34").Value = "  -4.66%  "
$ws.Range("E35").Value = "  -3.18%  "
$ws.Range("E36").Value = "  -3.40%  "
$ws.Range("E37").Value = "  -5.56%  "
$ws.Range("E38").Value = "  -6.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.32"
$ws.Range("E39").Value = "  -2.43%  "
$ws.Range("E40").Value = "  -6.60%  "
$ws.Range("E41").Value = "  -6.25%  "
$ws.Range("E42").Value = "  -6.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.292"
$ws.Range("E43").Value = "  -6.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.62"
$ws.Range("E44").Value = "  -9.60%  "
$ws.Range("E45").Value = "  -3.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "372.10"
$ws.Range("E46").Value = "  -5.69%  "
$ws.Range("D47").Value = "2.691.81"
$ws.Range("E47").Value = "  -1.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.24"
$ws.Range("E48").Value = "  -0.54%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("E50").Value = "  -7.27%  "
$ws.Range("E51").Value = "  -3.26%  "
